$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Build the "False" text value for H2 first (as a shared string), along with
# the desired style: Text number format (numFmtId 49) and an Arial 10pt
# black font. Using a formula + copy/paste-values round trip forces Excel to
# store the result as a genuine string cell (t="s") rather than
# re-interpreting "False" as a boolean the way a direct .Value assignment
# would.
$h2 = $ws.Range("H2")
$h2.Formula = '="False"'
$h2.Copy()
$h2.PasteSpecial(-4163)  # xlPasteValues
$h2.Font.Name = "Arial"
$h2.Font.Size = 10
$h2.Font.Color = 0
$h2.NumberFormat = "@"

# Propagate that exact value + style down to H3:H27 (replacing the boolean
# FALSE values that were there before). Doing this cell-by-cell (format then
# value) keeps each destination cell a real text cell instead of reverting
# back to boolean.
$src = $ws.Range("H2")
$src.Copy()
for ($row = 3; $row -le 27; $row++) {
    $dst = $ws.Range("H$row")
    $dst.PasteSpecial(-4122)  # xlPasteFormats
    $dst.PasteSpecial(-4163)  # xlPasteValues
}

# Extend the same style (but leave the cells empty) down to H28:H32.
for ($row = 28; $row -le 32; $row++) {
    $dst = $ws.Range("H$row")
    $dst.Font.Name = "Arial"
    $dst.Font.Size = 10
    $dst.Font.Color = 0
    $dst.NumberFormat = "@"
}

# Update the selection to match the author's last selection after extending
# the column (G28 active cell, G28:J33 selected).
$ws.Range("G28:J33").Select()

$wb.Saved = $false
